# Edit script: adjust categories (ajuste: corrigindo as categorias)
# - Adds two new columns: S "Idade ignorada", T "Total"
# - Adds two new rows: 7 "Outros", 8 "Total"
# - Fills in the new T column totals for existing rows 2-6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: add new trailing columns S1, T1 ---
$ws.Range("S1").Value = "Idade ignorada"
$ws.Range("T1").Value = "Total"

# --- Row 2: Doenças do aparelho circulatório ---
$ws.Range("S2").NumberFormat = "General"
$ws.Range("S2").Style = "Normal"
$ws.Range("T2").Value = 2078

# --- Row 3: Doenças do aparelho geniturinário ---
$ws.Range("S3").NumberFormat = "General"
$ws.Range("S3").Style = "Normal"
$ws.Range("T3").Value = 301

# --- Row 4: Doenças do aparelho respiratório ---
$ws.Range("S4").NumberFormat = "General"
$ws.Range("S4").Style = "Normal"
$ws.Range("T4").Value = 1084

# --- Row 5: Doenças endócrinas, nutricionais e metabólicas ---
$ws.Range("S5").NumberFormat = "General"
$ws.Range("S5").Style = "Normal"
$ws.Range("T5").Value = 240

# --- Row 6: Neoplasmas ---
$ws.Range("S6").NumberFormat = "General"
$ws.Range("S6").Style = "Normal"
$ws.Range("T6").Value = 1382

# --- Row 7 (new): Outros ---
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 137
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 28
$ws.Range("F7").Value = 80
$ws.Range("G7").Value = 71
$ws.Range("H7").Value = 75
$ws.Range("I7").Value = 94
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 82
$ws.Range("L7").Value = 118
$ws.Range("M7").Value = 112
$ws.Range("N7").Value = 115
$ws.Range("O7").Value = 123
$ws.Range("P7").Value = 122
$ws.Range("Q7").Value = 155
$ws.Range("R7").Value = 568
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 1973

# --- Row 8 (new): Total ---
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 157
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 13
$ws.Range("E8").Value = 39
$ws.Range("F8").Value = 98
$ws.Range("G8").Value = 91
$ws.Range("H8").Value = 113
$ws.Range("I8").Value = 169
$ws.Range("J8").Value = 173
$ws.Range("K8").Value = 239
$ws.Range("L8").Value = 361
$ws.Range("M8").Value = 478
$ws.Range("N8").Value = 535
$ws.Range("O8").Value = 635
$ws.Range("P8").Value = 700
$ws.Range("Q8").Value = 807
$ws.Range("R8").Value = 2443
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 7058
